$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 137; existing rows 137-219 shift down to 138-220
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new daily price record
$ws.Cells.Item(137, 1).Value = 9
$ws.Cells.Item(137, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(137, 3).Value = "Metropolitana"
$ws.Cells.Item(137, 4).Value = 44603
$ws.Cells.Item(137, 5).Value = 13
$ws.Cells.Item(137, 6).Value = 300000001
$ws.Cells.Item(137, 7).Value = "Rabanito"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 6400
$ws.Cells.Item(137, 11).Value = 3000
$ws.Cells.Item(137, 12).Value = 3500
$ws.Cells.Item(137, 13).Value = 3250
$ws.Cells.Item(137, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(137, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(137, 16).Value = 32
$ws.Cells.Item(137, 17).Value = 100
$ws.Cells.Item(137, 18).Value = "Hortaliza"
